# Add a new row for "model_V5" to the Matrix worksheet, mirroring the
# previous "model_V4.1" row (row 7), and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row index
$newRow = 8
$srcRow = 7

# Clone formatting (styles, number formats, fills) from the previous row,
# cell by cell, so only the columns actually populated on row 7 receive an
# explicit style on row 8 (column H stays completely empty, as on row 7).
$srcCols = @(1, 2, 3, 4, 5, 6, 7, 9, 10)
foreach ($col in $srcCols) {
    $ws.Cells.Item($srcRow, $col).Copy()
    $ws.Cells.Item($newRow, $col).PasteSpecial(-4122) # xlPasteFormats
}
$excel.CutCopyMode = 0

# Column A: shared string label for the new model version
$ws.Cells.Item($newRow, 1).Value = "model_V5"

# Columns B-D: raw counts (same values as the previous row, per the diff)
$ws.Cells.Item($newRow, 2).Value = 0
$ws.Cells.Item($newRow, 3).Value = 1011
$ws.Cells.Item($newRow, 4).Value = 989

# Column E: highlighted value cell
$ws.Cells.Item($newRow, 5).Value = 0

# Column F: total
$ws.Cells.Item($newRow, 6).Formula = "=SUM(B8:E8)"

# Column G: percent formula
$ws.Cells.Item($newRow, 7).Formula = "=E8/F8"

# Column I: combined count
$ws.Cells.Item($newRow, 9).Formula = "=E8+D8"

# Column J: percent formula
$ws.Cells.Item($newRow, 10).Formula = "=I8/F8"

# Update the active selection to E7, matching the post-edit workbook state
$ws.Range("E7").Select()
